# Added SPI Flash ROM
# This script edits the GD32F103RC BOM worksheet:
#  - C3/C4 capacitor changed from 10pF (C1634) to 12pF (C38523)
#  - C7..C18 (100n caps) designator list gains C18, qty 10 -> 11
#  - New U3 SOIC-8 W25Q128 SPI flash row inserted (row 16)
#  - J3 (USB) renamed/updated to J2 with new footprint/part (row 17)
#  - New Q1 8MHz crystal row (row 18)
#  - New Q2 32.768kHz crystal row (row 19)
#  - SW1/SW2 tactile switch row shifted down to row 20
#  - Grand total SUM formula relocated to row 22, now summing A2:A20 (42)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 : C3, C4 comment / LCSC part number change ----
$ws.Range("D3").Value = "12pF"
$ws.Range("E3").Value = "C38523"

# ---- Row 5 : C7..C17 group gains C18, quantity goes from 10 to 11 ----
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "C7, C8, C9, C10, C11, C12, C14, C15, C16, C17, C18"

# ---- Row 14 (U1) : E14 cell alignment loses its explicit left alignment ----
$ws.Range("E14").HorizontalAlignment = 1

# ---- Row 16 : brand new U3 SOIC-8 SPI flash component ----
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "U3"
$ws.Range("C16").Value = "SOIC-8_208mil"
$ws.Range("D16").Value = "W25Q128"
$ws.Range("D16").WrapText = $true
$ws.Range("E16").Value = "C97521"
$ws.Range("E16").WrapText = $true

# ---- Row 17 : USB connector re-designated J3 -> J2, new footprint/part ----
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "J2"
$ws.Range("C17").Value = "10118192-0002LF"
$ws.Range("C17").HorizontalAlignment = 1
$ws.Range("D17").Value = "USB-B-MICRO-SMD"
$ws.Range("E17").Value = "C2972784"
$ws.Range("E17").HorizontalAlignment = 1

# ---- Row 18 : new Q1 8 MHz crystal ----
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Q1"
$ws.Range("C18").Value = "HC-49S-SMD-2P"
$ws.Range("C18").HorizontalAlignment = 1
$ws.Range("D18").Value = "8 Mhz Crystal 20pF"
$ws.Range("E18").Value = "C12674"
$ws.Range("E18").HorizontalAlignment = 1

# ---- Row 19 : new Q2 32.768 kHz crystal (replaces old SUM row) ----
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Q2"
$ws.Range("C19").Value = "FC-12M"
$ws.Range("C19").HorizontalAlignment = 1
$ws.Range("D19").Value = "32.768kHz Crystal 12.5pF"
$ws.Range("D19").WrapText = $true
$ws.Range("E19").Value = "C32346"
$ws.Range("E19").WrapText = $true

# ---- Row 20 : SW1, SW2 tactile switches (moved down from row 17) ----
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "SW1, SW2"
$ws.Range("C20").Value = [char]0xFEFF + "EVQ-Q2"
$ws.Range("D20").Value = "SPST 6*6mm Tactile Switch"
$ws.Range("D20").WrapText = $true
$ws.Range("E20").Value = "C221880"

# ---- Row 22 : grand total, now summing rows 2 through 20 ----
$ws.Range("A22").Formula = "=SUM(A2:A20)"

# ---- Column E : width no longer flagged as a user customised width ----
$ws.Columns("E").ColumnWidth = 11.58

# ---- Selection / view state ----
$ws.Range("A1").Select()
$ws.Range("A22").Select()
